$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.162423666666667
$ws.Range("H2").Value = 3.487271
$ws.Range("I2").Value = 0.6447270069705344
$ws.Range("J2").Value = 0.6447270069705344
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 24.50599914239011
$ws.Range("R2").Value = 220.553992281511
$ws.Range("S2").Value = 0.03685298349584755
$ws.Range("T2").Value = 0.03685298349584755

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.162423666666667
$ws.Range("H3").Value = 3.487271
$ws.Range("I3").Value = 0.6447270069705344
$ws.Range("J3").Value = 0.6447270069705344
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 350.5871506803185
$ws.Range("R3").Value = 3155.284356122866
$ws.Range("S3").Value = 0.5272252889101288
$ws.Range("T3").Value = 0.5272252889101288

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.162423666666667
$ws.Range("H4").Value = 3.487271
$ws.Range("I4").Value = 0.6447270069705344
$ws.Range("J4").Value = 0.6447270069705344
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 53.62870608010879
$ws.Range("R4").Value = 482.6583547209791
$ws.Range("S4").Value = 0.0806487345645581
$ws.Range("T4").Value = 0.0806487345645581

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.4242653333333333
$ws.Range("H5").Value = 1.272796
$ws.Range("I5").Value = 0.2353146502133239
$ws.Range("J5").Value = 0.2353146502133239
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 8.944282702559555
$ws.Range("R5").Value = 80.498544323036
$ws.Range("S5").Value = 0.01345072693850887
$ws.Range("T5").Value = 0.01345072693850887

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.4242653333333333
$ws.Range("H6").Value = 1.272796
$ws.Range("I6").Value = 0.2353146502133239
$ws.Range("J6").Value = 0.2353146502133239
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 127.9584876074462
$ws.Range("R6").Value = 1151.626388467016
$ws.Range("S6").Value = 0.1924284745360072
$ws.Range("T6").Value = 0.1924284745360072

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.4242653333333333
$ws.Range("H7").Value = 1.272796
$ws.Range("I7").Value = 0.2353146502133239
$ws.Range("J7").Value = 0.2353146502133239
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 19.57358707824489
$ws.Range("R7").Value = 176.162283704204
$ws.Range("S7").Value = 0.02943544873880788
$ws.Range("T7").Value = 0.02943544873880788

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2162813333333333
$ws.Range("H8").Value = 0.648844
$ws.Range("I8").Value = 0.1199583428161417
$ws.Range("J8").Value = 0.1199583428161417
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 4.559602768911555
$ws.Range("R8").Value = 41.036424920204
$ws.Range("S8").Value = 0.006856891025498078
$ws.Range("T8").Value = 0.006856891025498078

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2162813333333333
$ws.Range("H9").Value = 0.648844
$ws.Range("I9").Value = 0.1199583428161417
$ws.Range("J9").Value = 0.1199583428161417
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 65.23048228715821
$ws.Range("R9").Value = 587.0743405844239
$ws.Range("S9").Value = 0.09809589371104327
$ws.Range("T9").Value = 0.09809589371104327

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2162813333333333
$ws.Range("H10").Value = 0.648844
$ws.Range("I10").Value = 0.1199583428161417
$ws.Range("J10").Value = 0.1199583428161417
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 9.97819331157289
$ws.Range("R10").Value = 89.80373980415601
$ws.Range("S10").Value = 0.01500555807960039
$ws.Range("T10").Value = 0.01500555807960039

Write-Output "done"